$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '62.847.53'
$ws.Range('E2').Value = '  +2.82%  '

# Row 3
$ws.Range('D3').Value = '2.964.95'
$ws.Range('E3').Value = '  +1.24%  '

# Row 4
$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  -0.12%  '

# Row 5
$ws.Range('D5').Value = '''594.22'
$ws.Range('E5').Value = '  +0.14%  '

# Row 6
$ws.Range('D6').Value = '''145.90'
$ws.Range('E6').Value = '  +0.46%  '

# Row 7
$ws.Range('E7').Value = '  -0.01%  '

# Row 8
$ws.Range('D8').Value = '2.962.57'
$ws.Range('E8').Value = '  +1.17%  '

# Row 9
$ws.Range('D9').Value = '''0.507'
$ws.Range('E9').Value = '  +0.54%  '

# Row 10
$ws.Range('D10').Value = '''7.25'
$ws.Range('E10').Value = '  +3.07%  '

# Row 11
$ws.Range('E11').Value = '  +1.28%  '

# Row 12
$ws.Range('D12').Value = '''0.444'
$ws.Range('E12').Value = '  +0.76%  '

# Row 13
$ws.Range('D13').Value = '''0.0000239'
$ws.Range('E13').Value = '  +5.95%  '

# Row 14
$ws.Range('D14').Value = '''33.32'
$ws.Range('E14').Value = '  -1.20%  '

# Row 15
$ws.Range('E15').Value = '  -0.24%  '

# Row 16
$ws.Range('D16').Value = '3.451.00'
$ws.Range('E16').Value = '  +0.97%  '

# Row 17
$ws.Range('D17').Value = '62.652.09'
$ws.Range('E17').Value = '  +2.43%  '

# Row 18
$ws.Range('B18').Value = 'Polkadot'
$ws.Range('C18').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D18').Value = '''6.71'
$ws.Range('E18').Value = '  -0.29%  '

# Row 19
$ws.Range('B19').Value = 'WrappedEther'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D19').Value = '2.957.73'
$ws.Range('E19').Value = '  +0.89%  '

# Row 20
$ws.Range('D20').Value = '''443.39'
$ws.Range('E20').Value = '  +1.97%  '

# Row 21
$ws.Range('D21').Value = '''13.44'
$ws.Range('E21').Value = '  -0.13%  '

# Row 22
$ws.Range('E22').Value = '  -1.12%  '

# Row 23
$ws.Range('E23').Value = '  -0.32%  '

# Row 24
$ws.Range('B24').Value = 'RenderToken'
$ws.Range('C24').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D24').Value = '''11.40'
$ws.Range('E24').Value = '  +3.18%  '

# Row 25
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').Value = '''81.88'
$ws.Range('E25').Value = '  +0.44%  '

# Row 26
$ws.Range('E26').Value = '  -2.86%  '

# Row 27
$ws.Range('E27').Value = '  +0.09%  '

# Row 28
$ws.Range('E28').Value = '  -0.07%  '

# Row 29
$ws.Range('D29').Value = '''7.17'
$ws.Range('E29').Value = '  +3.04%  '

# Row 30
$ws.Range('E30').Value = '  -0.23%  '

# Row 31
$ws.Range('D31').Value = '''2.13'
$ws.Range('E31').Value = '  -5.27%  '

# Row 32
$ws.Range('B32').Value = 'PEPE'
$ws.Range('C32').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D32').Value = '0.0₃0939'
$ws.Range('E32').Value = '  +8.18%  '

# Row 33
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').Value = '''26.70'
$ws.Range('E33').Value = '  -0.26%  '

# Row 34
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').Value = '''0.110'
$ws.Range('E34').Value = '  -0.18%  '

# Row 35
$ws.Range('E35').Value = '  -0.09%  '

# Row 36
$ws.Range('D36').Value = '''0.998'
$ws.Range('E36').Value = '  -1.41%  '

# Row 37
$ws.Range('E37').Value = '  -0.29%  '

# Row 38
$ws.Range('D38').Value = '''3.03'
$ws.Range('E38').Value = '  +1.41%  '

# Row 39
$ws.Range('D39').Value = '''2.05'
$ws.Range('E39').Value = '  +3.10%  '

# Row 40
$ws.Range('D40').Value = '''49.49'
$ws.Range('E40').Value = '  -0.93%  '

# Row 41
$ws.Range('D41').Value = '''8.55'
$ws.Range('E41').Value = '  -0.53%  '

# Row 42
$ws.Range('E42').Value = '  -4.08%  '

# Row 43
$ws.Range('E43').Value = '  -0.62%  '

# Row 44
$ws.Range('D44').Value = '''39.77'
$ws.Range('E44').Value = '  -5.78%  '

# Row 45
$ws.Range('D45').Value = '2.738.03'
$ws.Range('E45').Value = '  +1.13%  '

# Row 46
$ws.Range('D46').Value = '''136.74'
$ws.Range('E46').Value = '  +2.08%  '

# Row 47
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').Value = '''0.0340'
$ws.Range('E47').Value = '  -1.84%  '

# Row 48
$ws.Range('B48').Value = 'Bittensor'
$ws.Range('C48').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D48').Value = '''362.79'
$ws.Range('E48').Value = '  -2.82%  '

# Row 50
$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D50').Value = '''0.105'
$ws.Range('E50').Value = '  -0.33%  '

# Row 51
$ws.Range('B51').Value = 'InjectiveProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D51').Value = '''23.07'
$ws.Range('E51').Value = '  -3.14%  '
